$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (columns renamed to snake_case English names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish connector words (de/del/y/el/la/las/los) in municipality/state names ---
$ws.Range("B7").Value = 'Rincón De Romos'
$ws.Range("B8").Value = 'San Francisco De Los Romo'
$ws.Range("B9").Value = 'San José De Gracia'
$ws.Range("B14").Value = 'Playas De Rosarito'
$ws.Range("B30").Value = 'Amatenango De La Frontera'
$ws.Range("B34").Value = 'Benemérito De Las Américas'
$ws.Range("B41").Value = 'Chiapa De Corzo'
$ws.Range("B45").Value = 'Comitán De Domínguez'
$ws.Range("B67").Value = 'Mazapa De Madero'
$ws.Range("B72").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B82").Value = 'Salto De Agua'
$ws.Range("B83").Value = 'San Cristóbal De Las Casas'
$ws.Range("B126").Value = 'Guadalupe Y Calvo'
$ws.Range("B129").Value = 'Hidalgo Del Parral'
$ws.Range("B151").Value = 'San Francisco Del Oro'
$ws.Range("B169").Value = 'San Juan De Sabinas'
$ws.Range("B183").Value = 'Villa De Álvarez'
$ws.Range("A185").Value = 'Ciudad De México'
$ws.Range("B203").Value = 'Coneto De Comonfort'
$ws.Range("B216").Value = 'Nombre De Dios'
$ws.Range("B220").Value = 'Pánuco De Coronado'
$ws.Range("B227").Value = 'San Juan De Guadalupe'
$ws.Range("B228").Value = 'San Juan Del Río'
$ws.Range("B229").Value = 'San Luis Del Cordero'
$ws.Range("B230").Value = 'San Pedro Del Gallo'
$ws.Range("A240").Value = 'Estado De México'
$ws.Range("B240").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B243").Value = 'Almoloya De Alquisiras'
$ws.Range("B244").Value = 'Almoloya De Juárez'
$ws.Range("B248").Value = 'Atizapán De Zaragoza'
$ws.Range("B255").Value = 'Chapa De Mota'
$ws.Range("B258").Value = 'Coacalco De Berriozábal'
$ws.Range("B264").Value = 'Ecatepec De Morelos'
$ws.Range("B267").Value = 'Ixtapan De La Sal'
$ws.Range("B268").Value = 'Ixtapan Del Oro'
$ws.Range("B277").Value = 'Naucalpan De Juárez'
$ws.Range("B286").Value = 'San Felipe Del Progreso'
$ws.Range("B288").Value = 'San Simón De Guerrero'
$ws.Range("B297").Value = 'Tenango Del Valle'
$ws.Range("B303").Value = 'Tlalnepantla De Baz'
$ws.Range("B309").Value = 'Valle De Bravo'
$ws.Range("B310").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B311").Value = 'Villa De Allende'
$ws.Range("B312").Value = 'Villa Del Carbón'
$ws.Range("B325").Value = 'Apaseo El Alto'
$ws.Range("B333").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B337").Value = 'Jaral Del Progreso'
$ws.Range("B345").Value = 'Purísima Del Rincón'
$ws.Range("B349").Value = 'San Diego De La Unión'
$ws.Range("B351").Value = 'San Francisco Del Rincón'
$ws.Range("B353").Value = 'San Luis De La Paz'
$ws.Range("B355").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B357").Value = 'Silao De La Victoria'
$ws.Range("B361").Value = 'Valle De Santiago'
$ws.Range("B366").Value = 'Acapulco De Juárez'
$ws.Range("B368").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B369").Value = 'Alcozauca De Guerrero'
$ws.Range("B372").Value = 'Atenango Del Río'
$ws.Range("B374").Value = 'Atoyac De Álvarez'
$ws.Range("B375").Value = 'Ayutla De Los Libres'
$ws.Range("B378").Value = 'Buenavista De Cuéllar'
$ws.Range("B379").Value = 'Chilapa De Álvarez'
$ws.Range("B380").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B381").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B386").Value = 'Coyuca De Benítez'
$ws.Range("B387").Value = 'Coyuca De Catalán'
$ws.Range("B391").Value = 'Cuetzala Del Progreso'
$ws.Range("B392").Value = 'Cutzamala De Pinzón'
$ws.Range("B396").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B397").Value = 'Iguala De La Independencia'
$ws.Range("B399").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B400").Value = 'Zihuatanejo De Azueta'
$ws.Range("B402").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B417").Value = 'Taxco De Alarcón'
$ws.Range("B419").Value = 'Técpan De Galeana'
$ws.Range("B421").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B423").Value = 'Tixtla De Guerrero'
$ws.Range("B426").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B427").Value = 'Tlapa De Comonfort'
$ws.Range("B438").Value = 'Agua Blanca De Iturbide'
$ws.Range("B443").Value = 'Atotonilco El Grande'
$ws.Range("B448").Value = 'Cuautepec De Hinojosa'
$ws.Range("B453").Value = 'Huasca De Ocampo'
$ws.Range("B455").Value = 'Huejutla De Reyes'
$ws.Range("B461").Value = 'Mineral De La Reforma'
$ws.Range("B462").Value = 'Mineral Del Chico'
$ws.Range("B463").Value = 'Mineral Del Monte'
$ws.Range("B464").Value = 'Mixquiahuala De Juárez'
$ws.Range("B466").Value = 'Nopala De Villagrán'
$ws.Range("B467").Value = 'Omitlán De Juárez'
$ws.Range("B468").Value = 'Pachuca De Soto'
$ws.Range("B469").Value = 'Progreso De Obregón'
$ws.Range("B474").Value = 'Santiago De Anaya'
$ws.Range("B477").Value = 'Tenango De Doria'
$ws.Range("B479").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B481").Value = 'Tezontepec De Aldama'
$ws.Range("B487").Value = 'Tula De Allende'
$ws.Range("B488").Value = 'Tulancingo De Bravo'
$ws.Range("B491").Value = 'Zacualtipán De Ángeles'
$ws.Range("B492").Value = 'Zapotlán De Juárez'
$ws.Range("B496").Value = 'Acatlán De Juárez'
$ws.Range("B500").Value = 'Atemajac De Brizuela'
$ws.Range("B503").Value = 'Atotonilco El Alto'
$ws.Range("B505").Value = 'Autlán De Navarro'
$ws.Range("B515").Value = 'Cuautitlán De García Barragán'
$ws.Range("B519").Value = 'Encarnación De Díaz'
$ws.Range("B525").Value = 'Huejuquilla El Alto'
$ws.Range("B526").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B527").Value = 'Ixtlahuacán Del Río'
$ws.Range("B531").Value = 'Jilotlán De Los Dolores'
$ws.Range("B535").Value = 'La Manzanilla De La Paz'
$ws.Range("B536").Value = 'Lagos De Moreno'
$ws.Range("B543").Value = 'Ojuelos De Jalisco'
$ws.Range("B547").Value = 'San Cristóbal De La Barranca'
$ws.Range("B548").Value = 'San Diego De Alejandría'
$ws.Range("B551").Value = 'San Martín De Bolaños'
$ws.Range("B553").Value = 'San Miguel El Alto'
$ws.Range("B554").Value = 'Santa María De Los Ángeles'
$ws.Range("B557").Value = 'Talpa De Allende'
$ws.Range("B558").Value = 'Tamazula De Gordiano'
$ws.Range("B562").Value = 'Teocuitatlán De Corona'
$ws.Range("B563").Value = 'Tepatitlán De Morelos'
$ws.Range("B565").Value = 'Tizapán El Alto'
$ws.Range("B566").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B577").Value = 'Unión De San Antonio'
$ws.Range("B578").Value = 'Valle De Juárez'
$ws.Range("B583").Value = 'Yahualica De González Gallo'
$ws.Range("B584").Value = 'Zacoalco De Torres'
$ws.Range("B587").Value = 'Zapotitlán De Vadillo'
$ws.Range("B588").Value = 'Zapotlán Del Rey'
$ws.Range("B589").Value = 'Zapotlán El Grande'
$ws.Range("B609").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B611").Value = 'Cojumatlán De Régules'
$ws.Range("B705").Value = 'Puente De Ixtla'
$ws.Range("B709").Value = 'Tetela Del Volcán'
$ws.Range("B710").Value = 'Tlaltizapán De Zapata'
$ws.Range("B716").Value = 'Zacualpan De Amilpas'
$ws.Range("B720").Value = 'Amatlán De Cañas'
$ws.Range("B721").Value = 'Bahía De Banderas'
$ws.Range("B725").Value = 'Ixtlán Del Río'
$ws.Range("B731").Value = 'Santa María Del Oro'
$ws.Range("B746").Value = 'Montemorelos'
$ws.Range("B749").Value = 'San Nicolás De Los Garza'
$ws.Range("B751").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B757").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B758").Value = 'Coicoyán De Las Flores'
$ws.Range("B760").Value = 'Cuilápam De Guerrero'
$ws.Range("B761").Value = 'El Barrio De La Soledad'
$ws.Range("B762").Value = 'Guevea De Humboldt'
$ws.Range("B763").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B764").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B766").Value = 'Ixtlán De Juárez'
$ws.Range("B767").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B770").Value = 'Mariscala De Juárez'
$ws.Range("B772").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B773").Value = 'Nejapa De Madero'
$ws.Range("B774").Value = 'Oaxaca De Juárez'
$ws.Range("B775").Value = 'Ocotlán De Morelos'
$ws.Range("B776").Value = 'Putla Villa De Guerrero'
$ws.Range("B780").Value = 'San Agustín De Las Juntas'
$ws.Range("B820").Value = 'San Mateo Del Mar'
$ws.Range("B827").Value = 'San Miguel Del Puerto'
$ws.Range("B844").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B845").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B846").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B855").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B856").Value = 'Santa Inés De Zaragoza'
$ws.Range("B860").Value = 'Santa María Del Tule'
$ws.Range("B899").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B900").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B901").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B902").Value = 'Tlacolula De Matamoros'
$ws.Range("B904").Value = 'Villa De Etla'
$ws.Range("B905").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B906").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B907").Value = 'Villa De Zaachila'
$ws.Range("B908").Value = 'Villa Sola De Vega'
$ws.Range("B909").Value = 'Zimatlán De Álvarez'
$ws.Range("B925").Value = 'Chalchicomula De Sesma'
$ws.Range("B935").Value = 'Cuayuca De Andrade'
$ws.Range("B942").Value = 'Huehuetlán El Chico'
$ws.Range("B947").Value = 'Izúcar De Matamoros'
$ws.Range("B953").Value = 'Los Reyes De Juárez'
$ws.Range("B954").Value = 'Mazapiltepec De Juárez'
$ws.Range("B959").Value = 'Palmar De Bravo'
$ws.Range("B969").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B971").Value = 'San Salvador El Verde'
$ws.Range("B975").Value = 'Tecali De Herrera'
$ws.Range("B980").Value = 'Tepanco De López'
$ws.Range("B984").Value = 'Tepexi De Rodríguez'
$ws.Range("B986").Value = 'Tetela De Ocampo'
$ws.Range("B990").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B999").Value = 'Xayacatlán De Bravo'
$ws.Range("B1012").Value = 'Amealco De Bonfil'
$ws.Range("B1014").Value = 'Cadereyta De Montes'
$ws.Range("B1017").Value = 'Jalpan De Serra'
$ws.Range("B1018").Value = 'Landa De Matamoros'
$ws.Range("B1021").Value = 'Pinal De Amoles'
$ws.Range("B1024").Value = 'San Juan Del Río'
$ws.Range("B1035").Value = 'Ciudad Del Maíz'
$ws.Range("B1039").Value = 'Mexquitic De Carmona'
$ws.Range("B1047").Value = 'Santa María Del Río'
$ws.Range("B1049").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1052").Value = 'Villa De Arriaga'
$ws.Range("B1053").Value = 'Villa De Ramos'
$ws.Range("B1119").Value = 'Nacozari De García'
$ws.Range("B1132").Value = 'San Pedro De La Cueva'
$ws.Range("B1166").Value = 'Soto La Marina'
$ws.Range("B1172").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B1178").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1182").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1185").Value = 'San Pablo Del Monte'
$ws.Range("B1200").Value = 'Amatlán De Los Reyes'
$ws.Range("B1208").Value = 'Boca Del Río'
$ws.Range("B1209").Value = 'Camarón De Tejeda'
$ws.Range("B1213").Value = 'Castillo De Teayo'
$ws.Range("B1223").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1234").Value = 'Hueyapan De Ocampo'
$ws.Range("B1236").Value = 'Ixhuatlán De Madero'
$ws.Range("B1237").Value = 'Ixhuatlán Del Café'
$ws.Range("B1246").Value = 'Juchique De Ferrer'
$ws.Range("B1249").Value = 'Lerdo De Tejada'
$ws.Range("B1253").Value = 'Martínez De La Torre'
$ws.Range("B1255").Value = 'Medellín De Bravo'
$ws.Range("B1264").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1266").Value = 'Paso De Ovejas'
$ws.Range("B1268").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1276").Value = 'Sayula De Alemán'
$ws.Range("B1288").Value = 'Tlacotepec De Mejía'
$ws.Range("B1297").Value = 'Vega De Alatorre'
$ws.Range("B1314").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B1316").Value = 'Concepción Del Oro'
$ws.Range("B1318").Value = 'El Plateado De Joaquín Amaro'
$ws.Range("B1328").Value = 'Jiménez Del Teul'
$ws.Range("B1332").Value = 'Mezquital Del Oro'
$ws.Range("B1335").Value = 'Moyahua De Estrada'
$ws.Range("B1336").Value = 'Nochistlán De Mejía'
$ws.Range("B1337").Value = 'Noria De Ángeles'
$ws.Range("B1348").Value = 'Teúl De González Ortega'
$ws.Range("B1349").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1350").Value = 'Trinidad García De La Cadena'
$ws.Range("B1353").Value = 'Villa De Cos'

# --- Minor floating point precision touch-ups (recompute artifacts) ---
$ws.Range("D116").Value = 0.009453210010881391
$ws.Range("D134").Value = 0.009589227421109904

# --- Remove trailing footnote rows (1361-1366) that are no longer part of the data table ---
$ws.Range("A1361:D1366").EntireRow.Delete()
